# Auto-generated edit script: applies numeric corrections to several leve-profit
# calculation rows across all 8 sheets, matching refreshed market-price data.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1124.6296
$ws.Range("I28").Value = 1003.2778
$ws.Range("J28").Value = 1367.3334
$ws.Range("K28").Value = 1003.2778
$ws.Range("L28").Value = 1367.3334
$ws.Range("M28").Value = -518.2778
$ws.Range("N28").Value = -2337.3334
$ws.Range("H53").Value = 400.17648
$ws.Range("I53").Value = 466.63635
$ws.Range("K53").Value = 466.63635
$ws.Range("M53").Value = 170.36365
$ws.Range("H97").Value = 3269.2
$ws.Range("J97").Value = 3305.4167
$ws.Range("L97").Value = 9916.250100000001
$ws.Range("N97").Value = -10908.2501
$ws.Range("H116").Value = 9188.416999999999
$ws.Range("J116").Value = 6578
$ws.Range("L116").Value = 6578
$ws.Range("N116").Value = -13462
$ws.Range("H141").Value = 9075.147999999999
$ws.Range("I141").Value = 9959.583000000001
$ws.Range("K141").Value = 29878.749
$ws.Range("M141").Value = -24698.749

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 21834.934
$ws.Range("I28").Value = 5504.8
$ws.Range("K28").Value = 5504.8
$ws.Range("M28").Value = -5312.8
$ws.Range("H32").Value = 16028.924
$ws.Range("I32").Value = 15506.431
$ws.Range("K32").Value = 15506.431
$ws.Range("M32").Value = -15219.431
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H55").Value = 20011.25
$ws.Range("H99").Value = 21834.934
$ws.Range("I99").Value = 5504.8
$ws.Range("K99").Value = 5504.8
$ws.Range("M99").Value = -2509.8
$ws.Range("H102").Value = 2851.0908
$ws.Range("I102").Value = 1295.2858
$ws.Range("J102").Value = 5573.75
$ws.Range("K102").Value = 1295.2858
$ws.Range("L102").Value = 5573.75
$ws.Range("M102").Value = 326.7141999999999
$ws.Range("N102").Value = -8817.75
$ws.Range("H132").Value = 5342.3477
$ws.Range("I132").Value = 5465.8887
$ws.Range("J132").Value = 4897.6
$ws.Range("K132").Value = 16397.6661
$ws.Range("L132").Value = 14692.8
$ws.Range("M132").Value = -13867.6661
$ws.Range("N132").Value = -19752.8

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H111").Value = 149000
$ws.Range("J111").Value = 149000
$ws.Range("L111").Value = 149000
$ws.Range("N111").Value = -157180
$ws.Range("H122").Value = 110067.8
$ws.Range("J122").Value = 110067.8
$ws.Range("L122").Value = 110067.8
$ws.Range("N122").Value = -119867.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 178313.33
$ws.Range("J20").Value = 178313.33
$ws.Range("L20").Value = 178313.33
$ws.Range("N20").Value = -178785.33
$ws.Range("H22").Value = 5854.857
$ws.Range("I22").Value = 9232.615
$ws.Range("J22").Value = 366
$ws.Range("K22").Value = 9232.615
$ws.Range("L22").Value = 366
$ws.Range("M22").Value = -8882.615
$ws.Range("N22").Value = -1066
$ws.Range("H30").Value = 178313.33
$ws.Range("J30").Value = 178313.33
$ws.Range("L30").Value = 178313.33
$ws.Range("N30").Value = -178495.33
$ws.Range("H107").Value = 630.86957
$ws.Range("I107").Value = 496.3125
$ws.Range("K107").Value = 496.3125
$ws.Range("M107").Value = 1423.6875
$ws.Range("H116").Value = 99332
$ws.Range("J116").Value = 99332
$ws.Range("L116").Value = 99332
$ws.Range("N116").Value = -108510
$ws.Range("H120").Value = 15996.75
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 15996.75
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 15996.75
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -23254.75
$ws.Range("H121").Value = 51774.668
$ws.Range("J121").Value = 51774.668
$ws.Range("L121").Value = 51774.668
$ws.Range("N121").Value = -54394.668
$ws.Range("H128").Value = 178313.33
$ws.Range("J128").Value = 178313.33
$ws.Range("L128").Value = 178313.33
$ws.Range("N128").Value = -188273.33

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1578.675
$ws.Range("I5").Value = 1165.3448
$ws.Range("K5").Value = 3496.0344
$ws.Range("M5").Value = -3384.0344
$ws.Range("H7").Value = 217.73685
$ws.Range("I7").Value = 112.666664
$ws.Range("K7").Value = 337.999992
$ws.Range("M7").Value = -225.999992
$ws.Range("H26").Value = 567.13336
$ws.Range("I26").Value = 142.15
$ws.Range("J26").Value = 1417.1
$ws.Range("K26").Value = 426.45
$ws.Range("L26").Value = 4251.299999999999
$ws.Range("M26").Value = -138.45
$ws.Range("N26").Value = -4827.299999999999
$ws.Range("H34").Value = 816.125
$ws.Range("I34").Value = 645.8
$ws.Range("K34").Value = 1937.4
$ws.Range("M34").Value = -1853.4
$ws.Range("H60").Value = 272.625
$ws.Range("J60").Value = 302.5
$ws.Range("L60").Value = 907.5
$ws.Range("N60").Value = -1409.5
$ws.Range("H70").Value = 768.2
$ws.Range("I70").Value = 768.2
$ws.Range("K70").Value = 2304.6
$ws.Range("M70").Value = -1989.6
$ws.Range("H73").Value = 768.2
$ws.Range("I73").Value = 768.2
$ws.Range("K73").Value = 2304.6
$ws.Range("M73").Value = -1212.6
$ws.Range("H75").Value = 947.5
$ws.Range("J75").Value = 900
$ws.Range("L75").Value = 2700
$ws.Range("N75").Value = -4696
$ws.Range("H78").Value = 947.5
$ws.Range("J78").Value = 900
$ws.Range("L78").Value = 8100
$ws.Range("N78").Value = -18084
$ws.Range("H107").Value = 708.52
$ws.Range("I107").Value = 1322.1666
$ws.Range("J107").Value = 514.7368
$ws.Range("K107").Value = 3966.4998
$ws.Range("L107").Value = 1544.2104
$ws.Range("M107").Value = -2046.4998
$ws.Range("N107").Value = -5384.2104
$ws.Range("H109").Value = 19137.857
$ws.Range("I109").Value = 16996.6
$ws.Range("J109").Value = 24491
$ws.Range("K109").Value = 50989.8
$ws.Range("L109").Value = 73473
$ws.Range("M109").Value = -49949.8
$ws.Range("N109").Value = -75553
$ws.Range("H115").Value = 3181.5715
$ws.Range("I115").Value = 3181.5715
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 9544.7145
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -8369.7145
$ws.Range("N115").ClearContents()
$ws.Range("H135").Value = 1578.675
$ws.Range("I135").Value = 1165.3448
$ws.Range("K135").Value = 10488.1032
$ws.Range("M135").Value = -7953.103200000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 7500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7349
$ws.Range("N43").ClearContents()
$ws.Range("H80").Value = 11147474
$ws.Range("I80").Value = 16960772
$ws.Range("K80").Value = 16960772
$ws.Range("M80").Value = -16959774
$ws.Range("H83").Value = 11147474
$ws.Range("I83").Value = 16960772
$ws.Range("K83").Value = 84803860
$ws.Range("M83").Value = -84798868
$ws.Range("H109").Value = 127999
$ws.Range("J109").Value = 127999
$ws.Range("L109").Value = 127999
$ws.Range("N109").Value = -130079
$ws.Range("H132").Value = 5599.5
$ws.Range("I132").Value = 4323
$ws.Range("J132").Value = 7240.7144
$ws.Range("K132").Value = 12969
$ws.Range("L132").Value = 21722.1432
$ws.Range("M132").Value = -10439
$ws.Range("N132").Value = -26782.1432

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 330.44186
$ws.Range("I16").Value = 338.09756
$ws.Range("J16").Value = 173.5
$ws.Range("K16").Value = 338.09756
$ws.Range("L16").Value = 173.5
$ws.Range("M16").Value = -168.09756
$ws.Range("N16").Value = -513.5
$ws.Range("H22").Value = 1144.8667
$ws.Range("I22").Value = 1211.125
$ws.Range("J22").Value = 1069.1428
$ws.Range("K22").Value = 1211.125
$ws.Range("L22").Value = 1069.1428
$ws.Range("M22").Value = -916.125
$ws.Range("N22").Value = -1659.1428
$ws.Range("H27").Value = 1144.8667
$ws.Range("I27").Value = 1211.125
$ws.Range("J27").Value = 1069.1428
$ws.Range("K27").Value = 1211.125
$ws.Range("L27").Value = 1069.1428
$ws.Range("M27").Value = -1104.125
$ws.Range("N27").Value = -1283.1428
$ws.Range("H46").Value = 2469.3044
$ws.Range("I46").Value = 1001.4286
$ws.Range("J46").Value = 3111.5
$ws.Range("K46").Value = 1001.4286
$ws.Range("L46").Value = 3111.5
$ws.Range("M46").Value = -813.4286
$ws.Range("N46").Value = -3487.5
$ws.Range("H55").Value = 389.26086
$ws.Range("J55").Value = 453.0625
$ws.Range("L55").Value = 453.0625
$ws.Range("N55").Value = -799.0625
$ws.Range("H106").Value = 12400
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H108").Value = 89975.336
$ws.Range("J108").Value = 89975.336
$ws.Range("L108").Value = 89975.336
$ws.Range("N108").Value = -97655.336

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 174999
$ws.Range("J27").Value = 174999
$ws.Range("L27").Value = 174999
$ws.Range("N27").Value = -175137
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H111").Value = 59999.332
$ws.Range("J111").Value = 59999.332
$ws.Range("L111").Value = 59999.332
$ws.Range("N111").Value = -68179.33199999999
$ws.Range("H126").Value = 1413.7
$ws.Range("I126").Value = 1413.7
$ws.Range("K126").Value = 4241.1
$ws.Range("M126").Value = -1771.1
